# Auto-generated Excel COM-interop script
# Applies numeric updates to Shinryu_Profits workbook sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 6464616.5
$ws.Range("J19").Value = 9092028
$ws.Range("L19").Value = 9092028
$ws.Range("N19").Value = -9092378
$ws.Range("H28").Value = 398.4
$ws.Range("J28").Value = 939
$ws.Range("L28").Value = 939
$ws.Range("N28").Value = -1909
$ws.Range("H62").Value = 4393.769
$ws.Range("I62").Value = 3983.2856
$ws.Range("J62").Value = 4872.6665
$ws.Range("K62").Value = 3983.2856
$ws.Range("L62").Value = 4872.6665
$ws.Range("M62").Value = -3359.2856
$ws.Range("N62").Value = -6120.6665
$ws.Range("H65").Value = 4393.769
$ws.Range("I65").Value = 3983.2856
$ws.Range("J65").Value = 4872.6665
$ws.Range("K65").Value = 19916.428
$ws.Range("L65").Value = 24363.3325
$ws.Range("M65").Value = -16796.428
$ws.Range("N65").Value = -30603.3325
$ws.Range("H94").Value = 2441.6667
$ws.Range("I94").Value = 2441.6667
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2441.6667
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -1990.6667
$ws.Range("H116").Value = 2617.0833
$ws.Range("J116").Value = 2702
$ws.Range("L116").Value = 2702
$ws.Range("N116").Value = -9586
$ws.Range("H132").Value = 1988.3636
$ws.Range("I132").Value = 2052.1292
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 6156.3876
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -3626.3876
$ws.Range("N132").Value = -8060
$ws.Range("H137").Value = 24214.023
$ws.Range("I137").Value = 1386.3214
$ws.Range("J137").Value = 64162.5
$ws.Range("K137").Value = 4158.9642
$ws.Range("L137").Value = 192487.5
$ws.Range("M137").Value = -1608.9642
$ws.Range("N137").Value = -197587.5
$ws.Range("H138").Value = 2339.4912
$ws.Range("I138").Value = 1402.9722
$ws.Range("J138").Value = 3944.9524
$ws.Range("K138").Value = 4208.9166
$ws.Range("L138").Value = 11834.8572
$ws.Range("M138").Value = 931.0834000000004
$ws.Range("N138").Value = -22114.8572
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2628.375
$ws.Range("I2").Value = 3237
$ws.Range("J2").Value = 2263.2
$ws.Range("K2").Value = 3237
$ws.Range("L2").Value = 2263.2
$ws.Range("M2").Value = -3124
$ws.Range("N2").Value = -2489.2
$ws.Range("H4").Value = 468.2
$ws.Range("I4").Value = 460
$ws.Range("K4").Value = 460
$ws.Range("M4").Value = -344
$ws.Range("H74").Value = 4472.483
$ws.Range("I74").Value = 5503.5
$ws.Range("J74").Value = 1232.1428
$ws.Range("K74").Value = 5503.5
$ws.Range("L74").Value = 1232.1428
$ws.Range("M74").Value = -4629.5
$ws.Range("N74").Value = -2980.1428
$ws.Range("H77").Value = 4472.483
$ws.Range("I77").Value = 5503.5
$ws.Range("J77").Value = 1232.1428
$ws.Range("K77").Value = 27517.5
$ws.Range("L77").Value = 6160.714
$ws.Range("M77").Value = -23149.5
$ws.Range("N77").Value = -14896.714
$ws.Range("H97").Value = 534.9545000000001
$ws.Range("I97").Value = 512.0714
$ws.Range("J97").Value = 575
$ws.Range("K97").Value = 512.0714
$ws.Range("L97").Value = 575
$ws.Range("M97").Value = -16.07140000000004
$ws.Range("N97").Value = -1567
$ws.Range("H116").Value = 2628.375
$ws.Range("I116").Value = 3237
$ws.Range("J116").Value = 2263.2
$ws.Range("K116").Value = 3237
$ws.Range("L116").Value = 2263.2
$ws.Range("M116").Value = -943
$ws.Range("N116").Value = -6851.2
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2628.375
$ws.Range("I3").Value = 3237
$ws.Range("J3").Value = 2263.2
$ws.Range("K3").Value = 3237
$ws.Range("L3").Value = 2263.2
$ws.Range("M3").Value = -3123
$ws.Range("N3").Value = -2491.2
$ws.Range("H105").Value = 2586.98
$ws.Range("I105").Value = 1370.7727
$ws.Range("J105").Value = 2930.013
$ws.Range("K105").Value = 1370.7727
$ws.Range("L105").Value = 2930.013
$ws.Range("M105").Value = 376.2273
$ws.Range("N105").Value = -6424.013
$ws.Range("H107").Value = 2078.0232
$ws.Range("I107").Value = 2067.1943
$ws.Range("J107").Value = 2133.7144
$ws.Range("K107").Value = 2067.1943
$ws.Range("L107").Value = 2133.7144
$ws.Range("M107").Value = -147.1943000000001
$ws.Range("N107").Value = -5973.7144
$ws.Range("H134").Value = 1261.5358
$ws.Range("I134").Value = 1279.4509
$ws.Range("K134").Value = 3838.3527
$ws.Range("M134").Value = -1303.3527
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3484.26
$ws.Range("I31").Value = 2852.3684
$ws.Range("J31").Value = 5485.25
$ws.Range("K31").Value = 2852.3684
$ws.Range("L31").Value = 5485.25
$ws.Range("M31").Value = -2557.3684
$ws.Range("N31").Value = -6075.25
$ws.Range("H34").Value = 3484.26
$ws.Range("I34").Value = 2852.3684
$ws.Range("J34").Value = 5485.25
$ws.Range("K34").Value = 2852.3684
$ws.Range("L34").Value = 5485.25
$ws.Range("M34").Value = -2650.3684
$ws.Range("N34").Value = -5889.25
$ws.Range("H82").Value = 25000
$ws.Range("J82").Value = 25000
$ws.Range("L82").Value = 25000
$ws.Range("N82").Value = -25722
$ws.Range("H85").Value = 25000
$ws.Range("J85").Value = 25000
$ws.Range("L85").Value = 25000
$ws.Range("N85").Value = -27496
$ws.Range("H94").Value = 3572.182
$ws.Range("I94").Value = 1821
$ws.Range("K94").Value = 1821
$ws.Range("M94").Value = -1370
$ws.Range("H99").Value = 2260.0667
$ws.Range("I99").Value = 1850.7
$ws.Range("K99").Value = 1850.7
$ws.Range("M99").Value = -352.7
$ws.Range("H107").Value = 561.44116
$ws.Range("I107").Value = 500.7143
$ws.Range("J107").Value = 659.53845
$ws.Range("K107").Value = 500.7143
$ws.Range("L107").Value = 659.53845
$ws.Range("M107").Value = 1419.2857
$ws.Range("N107").Value = -4499.53845
$ws.Range("H126").Value = 2260.0667
$ws.Range("I126").Value = 1850.7
$ws.Range("K126").Value = 5552.1
$ws.Range("M126").Value = -3082.1
$ws.Range("H132").Value = 2528.4707
$ws.Range("I132").Value = 1726
$ws.Range("J132").Value = 3999.6667
$ws.Range("K132").Value = 5178
$ws.Range("L132").Value = 11999.0001
$ws.Range("M132").Value = -2648
$ws.Range("N132").Value = -17059.0001
$ws.Range("H134").Value = 3994.6553
$ws.Range("I134").Value = 2334.0527
$ws.Range("J134").Value = 7149.8
$ws.Range("K134").Value = 7002.158100000001
$ws.Range("L134").Value = 21449.4
$ws.Range("M134").Value = -4467.158100000001
$ws.Range("N134").Value = -26519.4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 770370.75
$ws.Range("I97").Value = 910210.9399999999
$ws.Range("J97").Value = 1250
$ws.Range("K97").Value = 910210.9399999999
$ws.Range("L97").Value = 1250
$ws.Range("M97").Value = -909714.9399999999
$ws.Range("N97").Value = -2242
$ws.Range("H107").Value = 280.5625
$ws.Range("I107").Value = 242.5
$ws.Range("J107").Value = 394.75
$ws.Range("K107").Value = 242.5
$ws.Range("L107").Value = 394.75
$ws.Range("M107").Value = 1677.5
$ws.Range("N107").Value = -4234.75
$ws.Range("H113").Value = 7171.2856
$ws.Range("I113").Value = 1994.3334
$ws.Range("J113").Value = 14073.889
$ws.Range("K113").Value = 1994.3334
$ws.Range("L113").Value = 14073.889
$ws.Range("M113").Value = 175.6666
$ws.Range("N113").Value = -18413.889
$ws.Range("H122").Value = 4169392.5
$ws.Range("I122").Value = 6668899
$ws.Range("J122").Value = 3548.111
$ws.Range("K122").Value = 20006697
$ws.Range("L122").Value = 10644.333
$ws.Range("M122").Value = -20004247
$ws.Range("N122").Value = -15544.333
$ws.Range("H126").Value = 1696.9333
$ws.Range("I126").Value = 1089
$ws.Range("J126").Value = 2391.7144
$ws.Range("K126").Value = 3267
$ws.Range("L126").Value = 7175.1432
$ws.Range("M126").Value = -797
$ws.Range("N126").Value = -12115.1432
$ws.Range("H132").Value = 2364.2778
$ws.Range("I132").Value = 2017.8928
$ws.Range("K132").Value = 6053.678400000001
$ws.Range("M132").Value = -3523.678400000001
$ws.Range("H141").Value = 69000
$ws.Range("J141").Value = 69000
$ws.Range("L141").Value = 69000
$ws.Range("N141").Value = -79360
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 556246.75
$ws.Range("I22").Value = 909427.25
$ws.Range("J22").Value = 1248.8572
$ws.Range("K22").Value = 909427.25
$ws.Range("L22").Value = 1248.8572
$ws.Range("M22").Value = -909132.25
$ws.Range("N22").Value = -1838.8572
$ws.Range("H27").Value = 556246.75
$ws.Range("I27").Value = 909427.25
$ws.Range("J27").Value = 1248.8572
$ws.Range("K27").Value = 909427.25
$ws.Range("L27").Value = 1248.8572
$ws.Range("M27").Value = -909320.25
$ws.Range("N27").Value = -1462.8572
$ws.Range("H40").Value = 4098.1763
$ws.Range("I40").Value = 3346.9
$ws.Range("J40").Value = 5171.4287
$ws.Range("K40").Value = 3346.9
$ws.Range("L40").Value = 5171.4287
$ws.Range("M40").Value = -3210.9
$ws.Range("N40").Value = -5443.4287
$ws.Range("H61").Value = 3860.4707
$ws.Range("I61").Value = 4668.885
$ws.Range("J61").Value = 1233.125
$ws.Range("K61").Value = 4668.885
$ws.Range("L61").Value = 1233.125
$ws.Range("M61").Value = -4466.885
$ws.Range("N61").Value = -1637.125
$ws.Range("H113").Value = 3860.4707
$ws.Range("I113").Value = 4668.885
$ws.Range("J113").Value = 1233.125
$ws.Range("K113").Value = 4668.885
$ws.Range("L113").Value = 1233.125
$ws.Range("M113").Value = -2498.885
$ws.Range("N113").Value = -5573.125
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2599.8462
$ws.Range("I96").Value = 1900
$ws.Range("J96").Value = 2658.1667
$ws.Range("K96").Value = 1900
$ws.Range("L96").Value = 2658.1667
$ws.Range("M96").Value = -527
$ws.Range("N96").Value = -5404.1667
$ws.Range("H104").Value = 29999.5
$ws.Range("J104").Value = 29999.5
$ws.Range("L104").Value = 29999.5
$ws.Range("N104").Value = -36987.5
$ws.Range("H107").Value = 654.625
$ws.Range("I107").Value = 654.4286
$ws.Range("J107").Value = 654.9
$ws.Range("K107").Value = 1963.2858
$ws.Range("L107").Value = 1964.7
$ws.Range("M107").Value = -43.28579999999988
$ws.Range("N107").Value = -5804.7
$ws.Range("H122").Value = 2176.0857
$ws.Range("I122").Value = 1606.9
$ws.Range("J122").Value = 2935
$ws.Range("K122").Value = 4820.700000000001
$ws.Range("L122").Value = 8805
$ws.Range("M122").Value = -2370.700000000001
$ws.Range("N122").Value = -13705
